# Rename several header cells across existing report sheets, then add a
# new "17_bi_percentiles_outliers_tota" sheet with percentile/outlier data.

$wb = $excel.ActiveWorkbook

# Sheet 1 (01_view_dataset_totals): avg_age -> avg_customer_age
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").Value = "avg_customer_age"

# Sheet 2 (02_view_monthly_transactions): total_qty -> total_qty_sold
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D1").Value = "total_qty_sold"

# Sheet 3 (24_bi_total_count_by_category_d): cnt -> order_cnt_by_cat,
# rows_with_category -> orders_by_category_totals
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B1").Value = "order_cnt_by_cat"
$ws3.Range("C1").Value = "orders_by_category_totals"

# Sheet 4 (03_view_product_category_perfor): total_qty -> total_qty_sold
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C1").Value = "total_qty_sold"

# Sheet 5 (18_bi_top_5_max_orders_by_total): total_amount -> top_orders_max_amount,
# total_count -> total_orders_same_amount
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A1").Value = "top_orders_max_amount"
$ws5.Range("B1").Value = "total_orders_same_amount"

# Sheet 6 (19_bi_top_5_min_orders_by_total): total_amount -> top_orders_min_amount,
# total_count -> total_orders_same_amount
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A1").Value = "top_orders_min_amount"
$ws6.Range("B1").Value = "total_orders_same_amount"

# New sheet: 17_bi_percentiles_outliers_tota (appended after the last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add($null, $lastSheet)
$ws7.Name = "17_bi_percentiles_outliers_tota"

$ws7.Range("A1").Value = "non_null_count"
$ws7.Range("B1").Value = "p25"
$ws7.Range("C1").Value = "p75"
$ws7.Range("D1").Value = "below_p25"
$ws7.Range("E1").Value = "between_p25_p75"
$ws7.Range("F1").Value = "above_p75"

$headerRange = $ws7.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws7.Range("A2").Value = 1000
$ws7.Range("B2").Value = 60
$ws7.Range("C2").Value = 900
$ws7.Range("D2").Value = 217
$ws7.Range("E2").Value = 581
$ws7.Range("F2").Value = 202
